$wb = $excel.ActiveWorkbook

$ipc = $wb.Worksheets.Item("IPC")
$trafo = $wb.Worksheets.Item("trafo")

# --- Clear the "trafo" sheet data block (rows 2:10), mirroring the
#     "clear trafo" part of the commit. Columns whose resulting style
#     matches the column default end up fully removed on save (Clear),
#     while columns with a distinct cell style / explicit blank style
#     are kept as empty cells (ClearContents).

# Columns A,B,C always fully cleared (style==column default everywhere)
$trafo.Range("A2:C10").Clear()

# Column D: rows 2:4 and 7:10 fully cleared, rows 5:6 keep their
# (non-default) style so only contents are cleared there.
$trafo.Range("D2:D4").Clear()
$trafo.Range("D7:D10").Clear()
$trafo.Range("D5:D6").ClearContents()

# Columns G,H always fully cleared
$trafo.Range("G2:H10").Clear()

# Column E: distinct (blank) style from column default -> keep cell, clear value
$trafo.Range("E2:E10").ClearContents()
$trafo.Range("E2:E10").Style = "Normal"

# Column F: rows 2:8 carry a distinct style -> keep cell, clear value.
# Rows 9:10 match the column default style -> fully cleared.
$trafo.Range("F2:F8").ClearContents()
$trafo.Range("F9:F10").Clear()

# Column I: distinct (blank) style from column default -> keep cell, clear value
$trafo.Range("I2:I10").ClearContents()
$trafo.Range("I2:I10").Style = "Normal"

# --- Selection / active-tab bookkeeping: move the "tab selected" /
#     active-window selection from IPC to trafo (the sheet the commit
#     message calls out), selecting the cleared data block.
$trafo.Range("A2:I10").Select()
$trafo.Activate()
